# fix import authentication error!
# Rename sheet "buff" -> "buffConfig", append 10 new buff rows (66-75) with
# their K/L text columns, and update the sheet view (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename the sheet -------------------------------------------------
$ws.Name = "buffConfig"

# --- 2. Append the new data rows (44-53) ---------------------------------
# Columns: A id, B ?, C id(copy), D, E, F bool, G, H, I, J, K text, L text/num

$rows = @(
    @{ Row=44; A=66; B=1; C=66; D=0; E=1; F=$true;  G=150; H=0;  I=1; J=0; K="AttributionType:11;ValueType:0;Value:200"; L="9,0,0" },
    @{ Row=45; A=67; B=1; C=67; D=1; E=1; F=$true;  G=90;  H=0;  I=1; J=4; K="AttributionType:1";                        L=0 },
    @{ Row=46; A=68; B=1; C=68; D=0; E=1; F=$false; G=1;   H=0;  I=1; J=0; K="AttributionType:30;ValueType:1;Value:100"; L=0 },
    @{ Row=47; A=69; B=1; C=69; D=0; E=1; F=$false; G=1;   H=0;  I=1; J=0; K="AttributionType:7;ValueType:0;Value:1000"; L=0 },
    @{ Row=48; A=70; B=1; C=70; D=0; E=1; F=$false; G=1;   H=0;  I=1; J=0; K="AttributionType:32;ValueType:2;Value:12";  L=0 },
    @{ Row=49; A=71; B=1; C=71; D=1; E=1; F=$true;  G=150; H=0;  I=1; J=0; K="AttributionType:40;ValueType:1;Value:150"; L=0 },
    @{ Row=50; A=72; B=1; C=72; D=1; E=1; F=$false; G=1;   H=0;  I=1; J=0; K="AttributionType:39;ValueType:2;Value:7";   L=0 },
    @{ Row=51; A=73; B=1; C=73; D=1; E=1; F=$true;  G=150; H=0;  I=1; J=0; K="AttributionType:12;ValueType:0;Value:-1000"; L=0 },
    @{ Row=52; A=74; B=1; C=74; D=1; E=1; F=$false; G=90;  H=90; I=1; J=0; K="AttributionType:39;ValueType:2;Value:8";   L=0 },
    @{ Row=53; A=75; B=1; C=75; D=1; E=1; F=$false; G=1;   H=0;  I=1; J=0; K="AttributionType:8;ValueType:0;Value:-100"; L=0 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $r.A
    $ws.Cells.Item($row, 2).Value  = $r.B
    $ws.Cells.Item($row, 3).Value  = $r.C
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
}

# --- 2b. Partial (mixed-font) runs inside K50 / K51, matching the source
#         workbook's rich-text spans on the trailing numeric value ------
$run1 = $ws.Cells.Item(50, 11).Characters(38, 1)
$run1.Font.Name = "宋体"
$run1.Font.Size = 11

$run2 = $ws.Cells.Item(51, 11).Characters(38, 5)
$run2.Font.Name = "宋体"
$run2.Font.Size = 11

# --- 3. Sheet view: zoom to 80% and move the active selection -----------
$excel.ActiveWindow.Zoom = 80
$ws.Range("L48").Select() | Out-Null
